# Generate Report for Handoff
# Update the status + handoff datetime for the "b07ffab4-5540-460d-9686-9f583923cf1a"
# file row across the Overview, zh-cn and de-de sheets to reflect that a new
# handoff report was generated for it.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-38-18 14:38:22"

# --- zh-cn sheet ---
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-18 14:38:19"

# --- de-de sheet ---
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("E3").Value = "2016-03-18 14:38:22"
